$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2,3,4) are cyclically rotated: row 3's original data
# moves into row 2, row 4's original data moves into row 3, and row 2's
# original data moves into row 4 (columns D, J, K, M, P only).

$ws.Range("D2").Value = 44804
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 9500
$ws.Range("M2").Value = 9750
$ws.Range("P2").Value = 542

$ws.Range("D3").Value = 44714
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = 9500
$ws.Range("P3").Value = 528

$ws.Range("D4").Value = 44792
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 9000
$ws.Range("M4").Value = 9500
$ws.Range("P4").Value = 528
